# Update cryptos list — refresh Price (D) and Volume(1h) % (E) values,
# and fix two mis-ordered rows (19/20: Polkadot/BitcoinCash, 37/38: Stacks/EnergySwap).
# Price cells are plain text (e.g. "66.629.99"), so we force a Text number
# format before writing, then restore the cell's original (unstyled) look,
# to stop Excel from auto-coercing numeric-looking strings into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.629.99'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.45%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.318.85'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.62%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.43'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.51%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('E9').Value = '  -3.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.403'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.894.31'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('E13').Value = '  -0.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.17'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '66.681.18'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.37%  '
$ws.Range('E16').Value = '  -2.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.327.54'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.64'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '432.68'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.50%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.521'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('E25').Value = '  -2.96%  '
$ws.Range('E26').Value = '  -0.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.48%  '
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('E29').Value = '  -2.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.79'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.32'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.46%  '
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.24'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.80'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.51'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.03'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.62%  '
$ws.Range('B37').Value = 'EnergySwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '27.25'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.84'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.818.52'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.790'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.44'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('E42').Value = '  -5.01%  '
$ws.Range('E43').Value = '  -2.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '24.40'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.36'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '323.82'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.47%  '
$ws.Range('E48').Value = '  -4.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.985'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.17'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.55%  '
$ws.Range('E51').Value = '  -1.47%  '
